$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: Id change
$ws.Range("A2").Value = 81392969

# P2: Lokalnamn change
$ws.Range("P2").Value = "550 m NV Örnanäs, Sk"

# S2: Noggrannhet change
$ws.Range("S2").Value = 50

# Y2: Startdatum change (keep as text, not an Excel date serial)
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2014-06-24"

# AA2: Slutdatum change (keep as text, not an Excel date serial)
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2014-06-24"

# AC2: Publik kommentar - remove content (cell cleared entirely)
$ws.Range("AC2").ClearContents()

# AR2: Samlings-nummer - remove content (cell cleared entirely)
$ws.Range("AR2").ClearContents()

# AW2: Rapportör change
$ws.Range("AW2").Value = "Charlotte Wigermo"

# AX2: Observatörer change
$ws.Range("AX2").Value = "Lars Åkerman"

# AY2: Projektnamn change
$ws.Range("AY2").Value = "Skånes Flora Millora 2008-2015"
